{"js": "// Replace the multiplication expressions in the document with the updated\n// values (the underlying numbers changed while the \"<3-digit>\u00d7<1-digit>=\"\n// pattern stayed the same). Each old expression is unique in the document,\n// so a plain text search + replace is unambiguous and keeps the existing\n// run formatting (font, size, etc.) intact.\nconst replacements = [\n  [\"616\u00d76=\", \"174\u00d73=\"],\n  [\"227\u00d72=\", \"142\u00d75=\"],\n  [\"977\u00d73=\", \"586\u00d73=\"],\n  [\"721\u00d77=\", \"820\u00d76=\"],\n  [\"658\u00d74=\", \"420\u00d76=\"],\n  [\"433\u00d75=\", \"439\u00d74=\"],\n  [\"348\u00d77=\", \"731\u00d72=\"],\n  [\"763\u00d73=\", \"120\u00d73=\"],\n  [\"220\u00d76=\", \"309\u00d79=\"],\n  [\"532\u00d77=\", \"309\u00d74=\"],\n  [\"205\u00d79=\", \"417\u00d74=\"],\n  [\"892\u00d73=\", \"967\u00d73=\"],\n  [\"850\u00d72=\", \"411\u00d78=\"],\n  [\"775\u00d74=\", \"621\u00d73=\"],\n  [\"119\u00d73=\", \"136\u00d72=\"],\n  [\"610\u00d75=\", \"366\u00d75=\"],\n  [\"500\u00d73=\", \"357\u00d73=\"],\n  [\"967\u00d78=\", \"412\u00d73=\"],\n  [\"175\u00d75=\", \"689\u00d76=\"],\n  [\"634\u00d76=\", \"439\u00d78=\"],\n  [\"683\u00d76=\", \"320\u00d73=\"],\n  [\"665\u00d76=\", \"598\u00d75=\"],\n  [\"232\u00d76=\", \"784\u00d77=\"],\n  [\"938\u00d78=\", \"823\u00d76=\"],\n  [\"201\u00d77=\", \"333\u00d77=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const result of results.items) {\n    result.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the multiplication expressions in the document with the updated\n# values (the underlying numbers changed while the \"<3-digit>x<1-digit>=\"\n# pattern stayed the same). Each old expression is unique in the document,\n# so Find/Replace on the whole-document range is unambiguous and the\n# existing run formatting (font, size, etc.) is left untouched.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"616\u00d76=\", \"174\u00d73=\"),\n    @(\"227\u00d72=\", \"142\u00d75=\"),\n    @(\"977\u00d73=\", \"586\u00d73=\"),\n    @(\"721\u00d77=\", \"820\u00d76=\"),\n    @(\"658\u00d74=\", \"420\u00d76=\"),\n    @(\"433\u00d75=\", \"439\u00d74=\"),\n    @(\"348\u00d77=\", \"731\u00d72=\"),\n    @(\"763\u00d73=\", \"120\u00d73=\"),\n    @(\"220\u00d76=\", \"309\u00d79=\"),\n    @(\"532\u00d77=\", \"309\u00d74=\"),\n    @(\"205\u00d79=\", \"417\u00d74=\"),\n    @(\"892\u00d73=\", \"967\u00d73=\"),\n    @(\"850\u00d72=\", \"411\u00d78=\"),\n    @(\"775\u00d74=\", \"621\u00d73=\"),\n    @(\"119\u00d73=\", \"136\u00d72=\"),\n    @(\"610\u00d75=\", \"366\u00d75=\"),\n    @(\"500\u00d73=\", \"357\u00d73=\"),\n    @(\"967\u00d78=\", \"412\u00d73=\"),\n    @(\"175\u00d75=\", \"689\u00d76=\"),\n    @(\"634\u00d76=\", \"439\u00d78=\"),\n    @(\"683\u00d76=\", \"320\u00d73=\"),\n    @(\"665\u00d76=\", \"598\u00d75=\"),\n    @(\"232\u00d76=\", \"784\u00d77=\"),\n    @(\"938\u00d78=\", \"823\u00d76=\"),\n    @(\"201\u00d77=\", \"333\u00d77=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2)\n}\n"}
